$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 16.71895933333333
$ws.Cells.Item(2, 8).Value = 50.156878
$ws.Cells.Item(2, 9).Value = 0.02912144738161902
$ws.Cells.Item(2, 10).Value = 0.03059269312988411
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.506715
$ws.Cells.Item(2, 14).Value = 1.520145
$ws.Cells.Item(2, 15).Value = 0.003122343715987576
$ws.Cells.Item(2, 16).Value = 0.003132472094339857
$ws.Cells.Item(2, 17).Value = 8.47174747859
$ws.Cells.Item(2, 18).Value = 76.24572730731001
$ws.Cells.Item(2, 19).Value = [double]"9.092716823246099E-05"
$ws.Cells.Item(2, 20).Value = [double]"9.583075752006463E-05"

# Row 3
$ws.Cells.Item(3, 7).Value = 16.71895933333333
$ws.Cells.Item(3, 8).Value = 50.156878
$ws.Cells.Item(3, 9).Value = 0.02912144738161902
$ws.Cells.Item(3, 10).Value = 0.03059269312988411
$ws.Cells.Item(3, 13).Value = 88.13219433333332
$ws.Cells.Item(3, 14).Value = 264.396583
$ws.Cells.Item(3, 15).Value = 0.5430646480820168
$ws.Cells.Item(3, 16).Value = 0.5448262620252092
$ws.Cells.Item(3, 17).Value = 1473.47857301643
$ws.Cells.Item(3, 18).Value = 13261.30715714787
$ws.Cells.Item(3, 19).Value = 0.0158148285739379
$ws.Cells.Item(3, 20).Value = 0.01666770264323906

# Row 4
$ws.Cells.Item(4, 7).Value = 16.71895933333333
$ws.Cells.Item(4, 8).Value = 50.156878
$ws.Cells.Item(4, 9).Value = 0.02912144738161902
$ws.Cells.Item(4, 10).Value = 0.03059269312988411
$ws.Cells.Item(4, 13).Value = 1.5741895
$ws.Cells.Item(4, 14).Value = 3.148379
$ws.Cells.Item(4, 15).Value = 0.009700049718478087
$ws.Cells.Item(4, 16).Value = 0.006487676741301404
$ws.Cells.Item(4, 17).Value = 26.31881023346033
$ws.Cells.Item(4, 18).Value = 157.912861400762
$ws.Cells.Item(4, 19).Value = 0.0002824794874757479
$ws.Cells.Item(4, 20).Value = 0.0001984755036725204

# Row 5
$ws.Cells.Item(5, 7).Value = 16.71895933333333
$ws.Cells.Item(5, 8).Value = 50.156878
$ws.Cells.Item(5, 9).Value = 0.02912144738161902
$ws.Cells.Item(5, 10).Value = 0.03059269312988411
$ws.Cells.Item(5, 13).Value = 72.07364666666666
$ws.Cells.Item(5, 14).Value = 216.22094
$ws.Cells.Item(5, 15).Value = 0.4441129584835175
$ws.Cells.Item(5, 16).Value = 0.4455535891391496
$ws.Cells.Item(5, 17).Value = 1204.996367625035
$ws.Cells.Item(5, 18).Value = 10844.96730862532
$ws.Cells.Item(5, 19).Value = 0.0129332121519729
$ws.Cells.Item(5, 20).Value = 0.01363068422545247

# Row 6
$ws.Cells.Item(6, 9).Value = 0.2708539632042961
$ws.Cells.Item(6, 10).Value = 0.2845377865576845
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 0.506715
$ws.Cells.Item(6, 14).Value = 1.520145
$ws.Cells.Item(6, 15).Value = 0.003122343715987576
$ws.Cells.Item(6, 16).Value = 0.003132472094339857
$ws.Cells.Item(6, 17).Value = 78.79437961213502
$ws.Cells.Item(6, 18).Value = 709.149416509215
$ws.Cells.Item(6, 19).Value = 0.0008456991699612639
$ws.Cells.Item(6, 20).Value = 0.0008913066761771772

# Row 7
$ws.Cells.Item(7, 9).Value = 0.2708539632042961
$ws.Cells.Item(7, 10).Value = 0.2845377865576845
$ws.Cells.Item(7, 13).Value = 88.13219433333332
$ws.Cells.Item(7, 14).Value = 264.396583
$ws.Cells.Item(7, 15).Value = 0.5430646480820168
$ws.Cells.Item(7, 16).Value = 0.5448262620252092
$ws.Cells.Item(7, 17).Value = 13704.59050225693
$ws.Cells.Item(7, 18).Value = 123341.3145203123
$ws.Cells.Item(7, 19).Value = 0.1470912122091606
$ws.Cells.Item(7, 20).Value = 0.1550236586551501

# Row 8
$ws.Cells.Item(8, 9).Value = 0.2708539632042961
$ws.Cells.Item(8, 10).Value = 0.2845377865576845
$ws.Cells.Item(8, 13).Value = 1.5741895
$ws.Cells.Item(8, 14).Value = 3.148379
$ws.Cells.Item(8, 15).Value = 0.009700049718478087
$ws.Cells.Item(8, 16).Value = 0.006487676741301404
$ws.Cells.Item(8, 17).Value = 244.7870796097155
$ws.Cells.Item(8, 18).Value = 1468.722477658293
$ws.Cells.Item(8, 19).Value = 0.002627296909528506
$ws.Cells.Item(8, 20).Value = 0.001845989179871673

# Row 9
$ws.Cells.Item(9, 9).Value = 0.2708539632042961
$ws.Cells.Item(9, 10).Value = 0.2845377865576845
$ws.Cells.Item(9, 13).Value = 72.07364666666666
$ws.Cells.Item(9, 14).Value = 216.22094
$ws.Cells.Item(9, 15).Value = 0.4441129584835175
$ws.Cells.Item(9, 16).Value = 0.4455535891391496
$ws.Cells.Item(9, 17).Value = 11207.48009331522
$ws.Cells.Item(9, 18).Value = 100867.320839837
$ws.Cells.Item(9, 19).Value = 0.1202897549156457
$ws.Cells.Item(9, 20).Value = 0.1267768320464856

# Row 10
$ws.Cells.Item(10, 7).Value = 194.8548433333333
$ws.Cells.Item(10, 8).Value = 584.56453
$ws.Cells.Item(10, 9).Value = 0.3394024086099587
$ws.Cells.Item(10, 10).Value = 0.3565493705749576
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.506715
$ws.Cells.Item(10, 14).Value = 1.520145
$ws.Cells.Item(10, 15).Value = 0.003122343715987576
$ws.Cells.Item(10, 16).Value = 0.003132472094339857
$ws.Cells.Item(10, 17).Value = 98.73587193965
$ws.Cells.Item(10, 18).Value = 888.6228474568501
$ws.Cells.Item(10, 19).Value = 0.001059730977714352
$ws.Cells.Item(10, 20).Value = 0.001116880953580495

# Row 11
$ws.Cells.Item(11, 7).Value = 194.8548433333333
$ws.Cells.Item(11, 8).Value = 584.56453
$ws.Cells.Item(11, 9).Value = 0.3394024086099587
$ws.Cells.Item(11, 10).Value = 0.3565493705749576
$ws.Cells.Item(11, 13).Value = 88.13219433333332
$ws.Cells.Item(11, 14).Value = 264.396583
$ws.Cells.Item(11, 15).Value = 0.5430646480820168
$ws.Cells.Item(11, 16).Value = 0.5448262620252092
$ws.Cells.Item(11, 17).Value = 17172.98491944455
$ws.Cells.Item(11, 18).Value = 154556.864275001
$ws.Cells.Item(11, 19).Value = 0.1843174495899561
$ws.Cells.Item(11, 20).Value = 0.1942574607977952

# Row 12
$ws.Cells.Item(12, 7).Value = 194.8548433333333
$ws.Cells.Item(12, 8).Value = 584.56453
$ws.Cells.Item(12, 9).Value = 0.3394024086099587
$ws.Cells.Item(12, 10).Value = 0.3565493705749576
$ws.Cells.Item(12, 13).Value = 1.5741895
$ws.Cells.Item(12, 14).Value = 3.148379
$ws.Cells.Item(12, 15).Value = 0.009700049718478087
$ws.Cells.Item(12, 16).Value = 0.006487676741301404
$ws.Cells.Item(12, 17).Value = 306.7384483994783
$ws.Cells.Item(12, 18).Value = 1840.43069039687
$ws.Cells.Item(12, 19).Value = 0.003292220238087815
$ws.Cells.Item(12, 20).Value = 0.002313177058604808

# Row 13
$ws.Cells.Item(13, 7).Value = 194.8548433333333
$ws.Cells.Item(13, 8).Value = 584.56453
$ws.Cells.Item(13, 9).Value = 0.3394024086099587
$ws.Cells.Item(13, 10).Value = 0.3565493705749576
$ws.Cells.Item(13, 13).Value = 72.07364666666666
$ws.Cells.Item(13, 14).Value = 216.22094
$ws.Cells.Item(13, 15).Value = 0.4441129584835175
$ws.Cells.Item(13, 16).Value = 0.4455535891391496
$ws.Cells.Item(13, 17).Value = 14043.89912969535
$ws.Cells.Item(13, 18).Value = 126395.0921672582
$ws.Cells.Item(13, 19).Value = 0.1507330078042004
$ws.Cells.Item(13, 20).Value = 0.158861851764977

# Row 14
$ws.Cells.Item(14, 7).Value = 82.82950199999999
$ws.Cells.Item(14, 8).Value = 165.659004
$ws.Cells.Item(14, 9).Value = 0.1442742299952585
$ws.Cells.Item(14, 10).Value = 0.1010420758958371
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.506715
$ws.Cells.Item(14, 14).Value = 1.520145
$ws.Cells.Item(14, 15).Value = 0.003122343715987576
$ws.Cells.Item(14, 16).Value = 0.003132472094339857
$ws.Cells.Item(14, 17).Value = 41.97095110593
$ws.Cells.Item(14, 18).Value = 251.82570663558
$ws.Cells.Item(14, 19).Value = 0.0004504737354046416
$ws.Cells.Item(14, 20).Value = 0.0003165114830978797

# Row 15
$ws.Cells.Item(15, 7).Value = 82.82950199999999
$ws.Cells.Item(15, 8).Value = 165.659004
$ws.Cells.Item(15, 9).Value = 0.1442742299952585
$ws.Cells.Item(15, 10).Value = 0.1010420758958371
$ws.Cells.Item(15, 13).Value = 88.13219433333332
$ws.Cells.Item(15, 14).Value = 264.396583
$ws.Cells.Item(15, 15).Value = 0.5430646480820168
$ws.Cells.Item(15, 16).Value = 0.5448262620252092
$ws.Cells.Item(15, 17).Value = 7299.94576679722
$ws.Cells.Item(15, 18).Value = 43799.67460078332
$ws.Cells.Item(15, 19).Value = 0.07835023393967901
$ws.Cells.Item(15, 20).Value = 0.05505037651759644

# Row 16
$ws.Cells.Item(16, 7).Value = 82.82950199999999
$ws.Cells.Item(16, 8).Value = 165.659004
$ws.Cells.Item(16, 9).Value = 0.1442742299952585
$ws.Cells.Item(16, 10).Value = 0.1010420758958371
$ws.Cells.Item(16, 13).Value = 1.5741895
$ws.Cells.Item(16, 14).Value = 3.148379
$ws.Cells.Item(16, 15).Value = 0.009700049718478087
$ws.Cells.Item(16, 16).Value = 0.006487676741301404
$ws.Cells.Item(16, 17).Value = 130.389332338629
$ws.Cells.Item(16, 18).Value = 521.5573293545159
$ws.Cells.Item(16, 19).Value = 0.00139946720404915
$ws.Cells.Item(16, 20).Value = 0.0006555283256822338

# Row 17
$ws.Cells.Item(17, 7).Value = 82.82950199999999
$ws.Cells.Item(17, 8).Value = 165.659004
$ws.Cells.Item(17, 9).Value = 0.1442742299952585
$ws.Cells.Item(17, 10).Value = 0.1010420758958371
$ws.Cells.Item(17, 13).Value = 72.07364666666666
$ws.Cells.Item(17, 14).Value = 216.22094
$ws.Cells.Item(17, 15).Value = 0.4441129584835175
$ws.Cells.Item(17, 16).Value = 0.4455535891391496
$ws.Cells.Item(17, 17).Value = 5969.824260723959
$ws.Cells.Item(17, 18).Value = 35818.94556434375
$ws.Cells.Item(17, 19).Value = 0.06407405511612568
$ws.Cells.Item(17, 20).Value = 0.04501965956946059

# Row 18
$ws.Cells.Item(18, 7).Value = 124.2078576666667
$ws.Cells.Item(18, 8).Value = 372.623573
$ws.Cells.Item(18, 9).Value = 0.2163479508088675
$ws.Cells.Item(18, 10).Value = 0.2272780738416368
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 0.506715
$ws.Cells.Item(18, 14).Value = 1.520145
$ws.Cells.Item(18, 15).Value = 0.003122343715987576
$ws.Cells.Item(18, 16).Value = 0.003132472094339857
$ws.Cells.Item(18, 17).Value = 62.937984597565
$ws.Cells.Item(18, 18).Value = 566.4418613780849
$ws.Cells.Item(18, 19).Value = 0.0006755126646748568
$ws.Cells.Item(18, 20).Value = 0.0007119422239642408

# Row 19
$ws.Cells.Item(19, 7).Value = 124.2078576666667
$ws.Cells.Item(19, 8).Value = 372.623573
$ws.Cells.Item(19, 9).Value = 0.2163479508088675
$ws.Cells.Item(19, 10).Value = 0.2272780738416368
$ws.Cells.Item(19, 13).Value = 88.13219433333332
$ws.Cells.Item(19, 14).Value = 264.396583
$ws.Cells.Item(19, 15).Value = 0.5430646480820168
$ws.Cells.Item(19, 16).Value = 0.5448262620252092
$ws.Cells.Item(19, 17).Value = 10946.71104960567
$ws.Cells.Item(19, 18).Value = 98520.39944645103
$ws.Cells.Item(19, 19).Value = 0.1174909237692831
$ws.Cells.Item(19, 20).Value = 0.1238270634114285

# Row 20
$ws.Cells.Item(20, 7).Value = 124.2078576666667
$ws.Cells.Item(20, 8).Value = 372.623573
$ws.Cells.Item(20, 9).Value = 0.2163479508088675
$ws.Cells.Item(20, 10).Value = 0.2272780738416368
$ws.Cells.Item(20, 13).Value = 1.5741895
$ws.Cells.Item(20, 14).Value = 3.148379
$ws.Cells.Item(20, 15).Value = 0.009700049718478087
$ws.Cells.Item(20, 16).Value = 0.006487676741301404
$ws.Cells.Item(20, 17).Value = 195.5267053563611
$ws.Cells.Item(20, 18).Value = 1173.160232138167
$ws.Cells.Item(20, 19).Value = 0.002098585879336866
$ws.Cells.Item(20, 20).Value = 0.00147450667347017

# Row 21
$ws.Cells.Item(21, 7).Value = 124.2078576666667
$ws.Cells.Item(21, 8).Value = 372.623573
$ws.Cells.Item(21, 9).Value = 0.2163479508088675
$ws.Cells.Item(21, 10).Value = 0.2845377865576845
$ws.Cells.Item(21, 13).Value = 72.07364666666666
$ws.Cells.Item(21, 14).Value = 216.22094
$ws.Cells.Item(21, 15).Value = 0.4441129584835175
$ws.Cells.Item(21, 17).Value = 8952.113246690957
$ws.Cells.Item(21, 18).Value = 80569.0192202186
$ws.Cells.Item(21, 19).Value = 0.09608292849557266
$ws.Cells.Item(21, 20).Value = 0.101264561532774
